# Update "想去人数" (F) and "最低票价" (G) figures on the 展览 and 全部类型
# sheets to reflect the latest scrape (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 sheet ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 137
$ws1.Range("G2").Value = 70
$ws1.Range("G3").Value = 70
$ws1.Range("F5").Value = 11457
$ws1.Range("F6").Value = 202
$ws1.Range("F9").Value = 11403
$ws1.Range("F10").Value = 464
$ws1.Range("F12").Value = 73
$ws1.Range("F14").Value = 5675
$ws1.Range("F16").Value = 3486
$ws1.Range("F17").Value = 177

# --- 全部类型 sheet ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 137
$ws4.Range("G2").Value = 70
$ws4.Range("G3").Value = 70
$ws4.Range("F7").Value = 11457
$ws4.Range("F8").Value = 202
$ws4.Range("F11").Value = 11403
$ws4.Range("F12").Value = 464
$ws4.Range("F14").Value = 73
$ws4.Range("F17").Value = 5675
$ws4.Range("F19").Value = 3486
$ws4.Range("F20").Value = 177
